# Edit slide 8 ("Answers (dont cheat)") of the presentation:
#  1. Split the run "src/Plasmodium.sh" (last paragraph of the content
#     placeholder) into two runs: "src/" and "Plasmodium.sh".
#  2. Append a new bulleted paragraph after it containing "git checkout
#     master", split into three runs: "git", " ", "checkout master".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# The content placeholder currently has 4 paragraphs; the 4th one is
# "edit/run src/Plasmodium.sh".
$para4 = $tr.Paragraphs(4, 1)
$para4Start = $para4.Start

# --- 1. Split "src/Plasmodium.sh" into "src/" + "Plasmodium.sh" ---
# "edit/run " is 9 characters, so "src/" begins at offset 9 within the
# paragraph. Re-assigning the same text to a sub-range forces a run
# split at that boundary without altering the visible text.
$srcRange = $tr.Characters($para4Start + 9, 4)
$srcRange.Text = "src/"

# --- 2. Add a new bulleted paragraph "git checkout master" ---
$cr = [char]13
$para4 = $tr.Paragraphs(4, 1)
$para4.InsertAfter($cr + "git checkout master") | Out-Null

$para5 = $tr.Paragraphs(5, 1)
$para5Start = $para5.Start

# Split into "git" | " " | "checkout master"
$gitRange = $tr.Characters($para5Start, 3)
$gitRange.Text = "git"

$spaceRange = $tr.Characters($para5Start + 3, 1)
$spaceRange.Text = " "

$restRange = $tr.Characters($para5Start + 4, 16)
$restRange.Text = "checkout master"
